$d = $word.ActiveDocument

$replacements = @(
    @('2024-07-02 Tuesday', '2024-07-03 Wednesday'),
    @('15×11=165', '63×79=4977'),
    @('96×37=3552', '87×45=3915'),
    @('56×98=5488', '46×12=552'),
    @('35×87=3045', '89×40=3560'),
    @('34×24=816', '38×90=3420'),
    @('74×44=3256', '83×31=2573'),
    @('44×70=3080', '53×72=3816'),
    @('82×90=7380', '36×66=2376'),
    @('40×81=3240', '76×26=1976'),
    @('11×43=473', '70×82=5740'),
    @('72×71=5112', '34×77=2618'),
    @('79×28=2212', '21×90=1890'),
    @('75×66=4950', '55×31=1705'),
    @('41×49=2009', '26×25=650'),
    @('90×75=6750', '94×91=8554'),
    @('49×67=3283', '78×75=5850'),
    @('64×96=6144', '39×28=1092'),
    @('88×60=5280', '58×92=5336'),
    @('60×30=1800', '37×92=3404'),
    @('79×61=4819', '41×93=3813'),
    @('60×74=4440', '96×70=6720'),
    @('52×51=2652', '95×78=7410'),
    @('81×67=5427', '83×15=1245'),
    @('86×27=2322', '13×48=624'),
    @('12×81=972', '51×98=4998'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
